# Aggiornamento fino a 1/09/2021
# Appends rows 358-366 to the sheet, continuing the daily COVID series.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use copy/paste-special (formats only) from the last existing data row (357)
# to row 358 so the new date cells inherit the same style (s="2", date format)
# as the rest of column A, without introducing any new style entries.
$ws.Range("A357").Copy()
$ws.Range("A358:A366").PasteSpecial(-4122)

$newRows = @(
    @(358, 44432, 0, 10, 62.61740763932373),
    @(359, 44433, 2, 12, 75.14088916718849),
    @(360, 44434, 0, 12, 75.14088916718849),
    @(361, 44435, 3, 11, 68.8791484032561),
    @(362, 44436, 2, 11, 68.8791484032561),
    @(363, 44437, 3, 14, 87.66437069505322),
    @(364, 44438, 6, 16, 100.187852222918),
    @(365, 44439, 0, 16, 100.187852222918),
    @(366, 44440, 1, 15, 93.9261114589856)
)

foreach ($r in $newRows) {
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value = $r[1]
    $ws.Cells.Item($row, 2).Value = $r[2]
    $ws.Cells.Item($row, 3).Value = $r[3]
    $ws.Cells.Item($row, 4).Value = $r[4]
}
